# Applies the "updated rmi files and update to 3.4.3" change:
#  - About sheet: add "Oregon" label (B1) and a date stamp (C1) for the region/version update.
#  - GbPbT sheet: break the formula link for VOC, CO, NOx, PM10, PM25, SOx, BC, OC rows
#    (columns B:C) and replace with literal 0 values (Oregon has no region-specific GWP
#    data for these pollutants yet), while leaving CO2, CH4, N2O, and F gases formulas intact.
#  - Make the GbPbT sheet the active/selected tab instead of About.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsGbPbT = $wb.Worksheets.Item("GbPbT")

# --- About sheet: record region name + update timestamp ---
$wsAbout.Range("B1").Value = "Oregon"
$wsAbout.Range("C1").Value = (Get-Date -Year 2022 -Month 10 -Day 27 -Hour 0 -Minute 0 -Second 0).Date
$wsAbout.Range("C1").NumberFormat = "m/d/yyyy"

# --- GbPbT sheet: zero out the region-specific pollutants, keep formulas for the rest ---
$zeroRows = @(3, 4, 5, 6, 7, 8, 9, 10)
foreach ($r in $zeroRows) {
    $wsGbPbT.Range("B$r").Value = 0
    $wsGbPbT.Range("C$r").Value = 0
}

# --- Make GbPbT the active/selected sheet (tabSelected moves from About to GbPbT) ---
$wsGbPbT.Activate()
$wsGbPbT.Range("F9").Select()
